# Apply cryptos-list refresh per Fri Mar  8 15:27:58 UTC 2024 GitHub Actions update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D (Price) and E (Volume(1h)) are stored as plain text in this sheet,
# even though many Price values look numeric (e.g. "43.01", "0.0000354").
# Force text format on a cell before writing so COM does not silently coerce
# the string into a number / strip formatting (trailing zeros, grouping dots, etc).
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# --- Per-row Price / Volume(1h) updates ---
Set-TextValue $ws.Range("D2") "68.857.64"
Set-TextValue $ws.Range("E2") "  +3.04%  "

Set-TextValue $ws.Range("D3") "3.956.47"
Set-TextValue $ws.Range("E3") "  +3.88%  "

Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.30%  "

Set-TextValue $ws.Range("D5") "479.58"
Set-TextValue $ws.Range("E5") "  +7.38%  "

Set-TextValue $ws.Range("D6") "149.76"
Set-TextValue $ws.Range("E6") "  +2.69%  "

Set-TextValue $ws.Range("D7") "0.625"
Set-TextValue $ws.Range("E7") "  +0.23%  "

Set-TextValue $ws.Range("D9") "0.730"
Set-TextValue $ws.Range("E9") "  -0.63%  "

Set-TextValue $ws.Range("D10") "0.169"
Set-TextValue $ws.Range("E10") "  +8.68%  "

Set-TextValue $ws.Range("D11") "0.0000354"
Set-TextValue $ws.Range("E11") "  +10.30%  "

Set-TextValue $ws.Range("D12") "43.01"
Set-TextValue $ws.Range("E12") "  -0.94%  "

Set-TextValue $ws.Range("D13") "4.592.21"
Set-TextValue $ws.Range("E13") "  +4.30%  "

Set-TextValue $ws.Range("D14") "10.41"
Set-TextValue $ws.Range("E14") "  +1.00%  "

Set-TextValue $ws.Range("D15") "14.75"
Set-TextValue $ws.Range("E15") "  -1.73%  "

Set-TextValue $ws.Range("D16") "3.967.91"
Set-TextValue $ws.Range("E16") "  +2.62%  "

Set-TextValue $ws.Range("E17") "  +0.11%  "

Set-TextValue $ws.Range("D18") "19.90"
Set-TextValue $ws.Range("E18") "  -0.11%  "

Set-TextValue $ws.Range("E19") "  -0.52%  "

Set-TextValue $ws.Range("D20") "68.896.17"
Set-TextValue $ws.Range("E20") "  +2.96%  "

Set-TextValue $ws.Range("D21") "438.41"
Set-TextValue $ws.Range("E21") "  +3.71%  "

Set-TextValue $ws.Range("D22") "3.32"
Set-TextValue $ws.Range("E22") "  +2.75%  "

Set-TextValue $ws.Range("D23") "14.39"
Set-TextValue $ws.Range("E23") "  -1.61%  "

Set-TextValue $ws.Range("D24") "87.80"
Set-TextValue $ws.Range("E24") "  +1.30%  "

Set-TextValue $ws.Range("D25") "3.68"
Set-TextValue $ws.Range("E25") "  +7.43%  "

Set-TextValue $ws.Range("D26") "38.43"
Set-TextValue $ws.Range("E26") "  +3.14%  "

Set-TextValue $ws.Range("D27") "9.84"
Set-TextValue $ws.Range("E27") "  +3.97%  "

Set-TextValue $ws.Range("D28") "10.01"
Set-TextValue $ws.Range("E28") "  +3.56%  "

Set-TextValue $ws.Range("D29") "728.90"
Set-TextValue $ws.Range("E29") "  -2.45%  "

Set-TextValue $ws.Range("D30") "13.17"
Set-TextValue $ws.Range("E30") "  -3.60%  "

Set-TextValue $ws.Range("D31") "0.126"
Set-TextValue $ws.Range("E31") "  -4.96%  "

Set-TextValue $ws.Range("E32") "  +3.67%  "

Set-TextValue $ws.Range("D33") "42.26"
Set-TextValue $ws.Range("E33") "  -1.89%  "

Set-TextValue $ws.Range("D34") "60.17"
Set-TextValue $ws.Range("E34") "  +2.61%  "

Set-TextValue $ws.Range("D35") "0.0₃0845"
Set-TextValue $ws.Range("E35") "  +25.41%  "

Set-TextValue $ws.Range("E36") "  -2.51%  "

Set-TextValue $ws.Range("E37") "  -0.05%  "

Set-TextValue $ws.Range("D38") "5.38"
Set-TextValue $ws.Range("E38") "  -1.64%  "

Set-TextValue $ws.Range("D39") "0.0471"
Set-TextValue $ws.Range("E39") "  -0.63%  "

Set-TextValue $ws.Range("D40") "3.04"
Set-TextValue $ws.Range("E40") "  +5.93%  "

Set-TextValue $ws.Range("D41") "2.94"
Set-TextValue $ws.Range("E41") "  +10.80%  "

Set-TextValue $ws.Range("D42") "2.60"
Set-TextValue $ws.Range("E42") "  +5.09%  "

Set-TextValue $ws.Range("D43") "0.141"
Set-TextValue $ws.Range("E43") "  +0.75%  "

Set-TextValue $ws.Range("D46") "2.17"
Set-TextValue $ws.Range("E46") "  +2.99%  "

Set-TextValue $ws.Range("D48") "148.92"
Set-TextValue $ws.Range("E48") "  +1.21%  "

Set-TextValue $ws.Range("D50") "2.90"
Set-TextValue $ws.Range("E50") "  +1.65%  "

Set-TextValue $ws.Range("D51") "24.76"
Set-TextValue $ws.Range("E51") "  -1.40%  "

# --- Row 44/45: coin order swapped (TheGraph <-> FirstDigitalUSD) ---
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue $ws.Range("D44") "1.00"
Set-TextValue $ws.Range("E44") "  -0.02%  "

$ws.Range("B45").Value = "TheGraph"
$ws.Range("C45").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
Set-TextValue $ws.Range("D45") "0.333"
Set-TextValue $ws.Range("E45") "  -4.53%  "

# --- Row 47/49: coin order swapped (LidoDAOToken <-> ApeXProtocol); row 48 (Monero) unaffected ---
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue $ws.Range("D47") "3.23"
Set-TextValue $ws.Range("E47") "  -1.03%  "

$ws.Range("B49").Value = "LidoDAOToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D49") "3.38"
Set-TextValue $ws.Range("E49") "  -0.54%  "
